$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Numeric data refresh (Casos totales / Nuevos casos / Casos activos / Recuperados / Muertes hoy / Muertes) ---

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 1593303
$ws.Range("C4").Value = 580
$ws.Range("D4").Value = 370971
$ws.Range("E4").Value = 1127384

# Row 23 - Paises Bajos
$ws.Range("B23").Value = 44700
$ws.Range("C23").Value = 253
$ws.Range("G23").Value = 27
$ws.Range("H23").Value = 5775

# Row 27 - Suecia
$ws.Range("B27").Value = 32172
$ws.Range("C27").Value = 649
$ws.Range("E27").Value = 23330
$ws.Range("G27").Value = 40
$ws.Range("H27").Value = 3871

# Row 83 - Croacia
$ws.Range("B83").Value = 2237
$ws.Range("C83").Value = 3
$ws.Range("E83").Value = 162
$ws.Range("G83").Value = 1
$ws.Range("H83").Value = 97

# Row 108 - Libano
$ws.Range("D108").Value = 663
$ws.Range("E108").Value = 335

# Row 150 - Liberia
$ws.Range("B150").Value = 240
$ws.Range("C150").Value = 2
$ws.Range("D150").Value = 131
$ws.Range("E150").Value = 86

# Rows 199/200 - low countries re-sorted: swap underlying data so the
# country names (which keep their row position) end up with the correct stats.
# Row 199 was "Santa Lucia" (18,0,18,0,0,0,0) -> becomes "Belice" data (18,0,16,0,0,0,2)
# Row 200 was "Belice" (18,0,16,0,0,0,2) -> becomes "Santa Lucia" data (18,0,18,0,0,0,0)
$ws.Range("A199").Value = "Belice"
$ws.Range("D199").Value = 16
$ws.Range("H199").Value = 2

$ws.Range("A200").Value = "Santa Lucia"
$ws.Range("D200").Value = 18
$ws.Range("H200").Value = 0

# Rows 209/210 - Seychelles / Groenlandia swap places (tied totals, values identical)
$ws.Range("A209").Value = "Groenlandia"
$ws.Range("A210").Value = "Seychelles"

# Rows 214/215 - Bonaire, San Eustaquio y Saba / Sahara Occidental swap places (tied totals, values identical)
$ws.Range("A214").Value = "Sahara Occidental"
$ws.Range("A215").Value = "Bonaire, San Eustaquio y Saba"

# --- Update "last refreshed" timestamp string ---
$ws.Range("A1").Value = "Datos actualizados a 21 de Mayo de 2020 a las 14:35"
